# Daily attendance processing - 2025-12-13 07:50:08
# Normalizes the "Recorded By" (column G) entries so that the first two
# comma-separated recorders are swapped for a known set of values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping of old "Recorded By" text -> new text (first two entries swapped).
$map = @{
    "System, backup@backdoor.com, system" = "backup@backdoor.com, System, system"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
